$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks so they can be rebuilt consistently with the new URLs
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2026-02-04 12:54:41"
$ws.Range("B2").Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value = "2026-02-04 12:54:41"
$ws.Range("B3").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Range("A4").Value = "2026-02-04 12:54:41"
$ws.Range("B4").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

# Row 5
$ws.Range("A5").Value = "2026-02-04 12:54:41"
$ws.Range("B5").Value = "【Python】特定サイトからのデータ収集ツール開発依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5485630"
$ws.Range("G5").Value = 340
$ws.Range("H5").Value = "🔥Python ◆ツール,開発 ◇サイト"

# Row 6
$ws.Range("A6").Value = "2026-02-04 12:54:41"
$ws.Range("B6").Value = "【急募】Telegramグループ運営用BOT開発のフリーランス募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5485147"
$ws.Range("G6").Value = 180
$ws.Range("H6").Value = "★bot ◆開発"

# Row 7
$ws.Range("A7").Value = "2026-02-04 12:54:41"
$ws.Range("B7").Value = "防火防災管理システムの開発"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5485816"
$ws.Range("G7").Value = 108
$ws.Range("H7").Value = "◆開発 ◇管理"

# Row 8
$ws.Range("A8").Value = "2026-02-04 12:54:41"
$ws.Range("B8").Value = "【募集】会社Webサイト改善(WordPress/スターサーバー/小規模改修/GA4導入)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5485476"
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = "◇サイト ○WordPress"

# Row 9
$ws.Range("A9").Value = "2026-02-04 12:54:41"
$ws.Range("B9").Value = "ファイルメーカーでの在庫・顧客管理システム構築"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5485054"
$ws.Range("G9").Value = 53
$ws.Range("H9").Value = "◇管理"

# Row 10
$ws.Range("A10").Value = "2026-02-04 12:54:41"
$ws.Range("B10").Value = "bubbleで構築したサイトの修正対応"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5485362"
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "◇サイト"

# Row 11
$ws.Range("A11").Value = "2026-02-04 12:54:41"
$ws.Range("B11").Value = "【急募】日本人エンジニアを渋谷・横浜で緊急募集!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5484588"
$ws.Range("G11").Value = 25

# Row 12
$ws.Range("A12").Value = "2026-02-04 12:54:41"
$ws.Range("B12").Value = "【美容室向け】社内動画学習プラットフォーム構築依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5485460"
$ws.Range("G12").Value = 18

# Row 13
$ws.Range("A13").Value = "2026-02-04 12:54:41"
$ws.Range("B13").Value = "美容医療機器エンジニア募集!【HIFU・RF・CO2レーザーの機序を解明する専門家・プロを募集】"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5485304"
$ws.Range("G13").Value = 18

# Row 14
$ws.Range("A14").Value = "2026-02-04 12:54:41"
$ws.Range("B14").Value = "【継続依頼あり/高単価】Yoom/kintone 実装パートナー募集!"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5485174"
$ws.Range("G14").Value = 18

# Row 15
$ws.Range("A15").Value = "2026-02-04 12:54:41"
$ws.Range("B15").Value = "AWS環境の運用まわりのご相談"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5485814"
$ws.Range("G15").Value = 13

# Row 16
$ws.Range("A16").Value = "2026-02-04 12:54:41"
$ws.Range("B16").Value = "AKASHI初期設定・マネーフォワード給与連携"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5485734"
$ws.Range("G16").Value = 13

# Row 17
$ws.Range("A17").Value = "2026-02-04 12:54:41"
$ws.Range("B17").Value = "【急募】FX自動売買(.mp4)のエラー修正をお手伝いください!"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5485785"
$ws.Range("G17").Value = 10

# Row 18
$ws.Range("A18").Value = "2026-02-04 12:54:41"
$ws.Range("B18").Value = "Claude Consoleにプロンプトを入れてください"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5485506"
$ws.Range("G18").Value = 10

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5455098") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445159") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445154") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5485630") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5485147") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5485816") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5485476") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5485054") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5485362") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5484588") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5485460") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5485304") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5485174") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5485814") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5485734") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5485785") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5485506") | Out-Null

$ws.Range("F2:F18").Style = "Hyperlink"

$ws.Columns.Item(8).ColumnWidth = 21.166666666666668
